$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the attribute rows (A2:B7) into their new sequence.
$ws.Range("A2").Value = "SubProcessID"
$ws.Range("B2").Value = "str"

$ws.Range("A3").Value = "time:timestamp"
$ws.Range("B3").Value = "datetime"

$ws.Range("A4").Value = "stream:datastream"
$ws.Range("B4").Value = "dict"

$ws.Range("A5").Value = "operation_end_time"
$ws.Range("B5").Value = "datetime"

$ws.Range("A6").Value = "org:resource"
$ws.Range("B6").Value = "str"

$ws.Range("A7").Value = "concept:name"
$ws.Range("B7").Value = "str"

# Clear the bold/bordered/centered header style from the header row so it
# falls back to the default (unstyled) cell format.
$ws.Range("A1:B1").ClearFormats()
